$d = $word.ActiveDocument

# 1) "Curso (semestre ideal): EM (10)" -> "Curso (semestre ideal): EM (8)"
$d.Content.Find.Execute(
    "Curso (semestre ideal): EM (10)", $true, $false, $false, $false, $false,
    $true, 1, $false, "Curso (semestre ideal): EM (8)", 2) | Out-Null

# 2) Append a new "Requisitos" section (Heading2 + ListBullet entry) after
#    the Bibliografia paragraph, i.e. at the very end of the document body.
$wordNs = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"

$lastPara = $d.Paragraphs.Last
$lastPara.Range.InsertParagraphAfter() | Out-Null

$headingPara = $d.Paragraphs.Last
$headingPara.Range.InsertXML(
    "<w:p $wordNs><w:pPr><w:pStyle w:val='Heading2'/></w:pPr><w:r><w:t>Requisitos</w:t></w:r></w:p>")

$headingPara = $d.Paragraphs.Last
$headingPara.Range.InsertParagraphAfter() | Out-Null

$reqPara = $d.Paragraphs.Last
$reqPara.Range.InsertXML(
    "<w:p $wordNs><w:pPr><w:pStyle w:val='ListBullet'/></w:pPr><w:r><w:t>LOM3013 -  Ciência dos Materiais  (Requisito fraco)</w:t><w:br/></w:r></w:p>")
